# Week 15 simulations: add a new player column "E.Winston" (receiver),
# inserted before the existing "A.Trautman" column on both the "Rushing"
# and "Receiving" sheets. Inserting the column shifts A.Trautman,
# J.Johnson, G.Griffin and N.Vannett one column to the right and extends
# the used range from column Y to column Z.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    # Column V (22) currently holds "A.Trautman" - insert a new blank
    # column there, pushing A.Trautman..N.Vannett right by one (V:Y -> W:Z).
    $ws.Columns("V").Insert()

    # Populate the newly inserted column: header name in row 1, and the
    # same "n" placeholder used by every other player column in row 2.
    $ws.Cells.Item(1, 22).Value = "E.Winston"
    $ws.Cells.Item(2, 22).Value = "n"
}
